# "Added bal checks to picture ppt"
# Reposition/resize the four scanned-check pictures on slide 7 so they are
# evenly arranged (2x2 grid) instead of their original ad-hoc placement.
#
# PowerPoint's COM object model reports/accepts Shape.Left/Top/Width/Height
# in points, while the underlying OOXML stores EMUs (1 pt = 12700 EMU).
# The literal point values below were chosen so that, after PowerPoint's
# internal (single-precision) point<->EMU conversion, they land exactly on
# the target EMU coordinates.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# Picture 2 -> top-left
$shape1 = $s.Shapes.Item(1)
$shape1.Left   = 15.75007874015748
$shape1.Top    = 2.0
$shape1.Width  = 284.842125984252
$shape1.Height = 237.28669291338582

# Picture 3 -> top-right
$shape2 = $s.Shapes.Item(2)
$shape2.Left   = 303.9685039370079
$shape2.Top    = 2.0
$shape2.Width  = 285.218188976378
$shape2.Height = 237.6

# Picture 4 -> bottom-left
$shape3 = $s.Shapes.Item(3)
$shape3.Left   = 18.75023622047244
$shape3.Top    = 252.49503937007873
$shape3.Width  = 285.2182769775391
$shape3.Height = 237.6

# Picture 5 -> bottom-right
$shape4 = $s.Shapes.Item(4)
$shape4.Left   = 303.96842956542974
$shape4.Top    = 253.49503937007873
$shape4.Width  = 285.218188976378
$shape4.Height = 237.6
